# Update the date line and all the three-digit x one-digit multiplication
# problems to the new set of values described by the commit diff.

$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-09-28 Sunday"; new = "2025-09-29 Monday"},
    @{old = "534×3="; new = "175×6="},
    @{old = "234×2="; new = "749×2="},
    @{old = "648×9="; new = "315×3="},
    @{old = "169×7="; new = "477×6="},
    @{old = "604×7="; new = "423×2="},
    @{old = "889×9="; new = "890×8="},
    @{old = "696×7="; new = "956×4="},
    @{old = "193×4="; new = "161×2="},
    @{old = "687×7="; new = "373×5="},
    @{old = "800×9="; new = "736×3="},
    @{old = "974×6="; new = "860×5="},
    @{old = "508×2="; new = "327×7="},
    @{old = "965×4="; new = "938×4="},
    @{old = "376×9="; new = "713×3="},
    @{old = "705×8="; new = "451×6="},
    @{old = "105×9="; new = "811×4="},
    @{old = "195×8="; new = "707×2="},
    @{old = "901×2="; new = "569×9="},
    @{old = "289×6="; new = "206×6="},
    @{old = "966×6="; new = "867×7="},
    @{old = "780×8="; new = "151×5="},
    @{old = "399×9="; new = "324×8="},
    @{old = "227×2="; new = "765×4="},
    @{old = "284×4="; new = "733×6="},
    @{old = "549×3="; new = "623×9="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
